# ReadMe.docx small fix
#
# The only substantive (visible-text) change in the target revision is in
# the "Testing Notes" section: the paragraph that used to read
#
#   "Finally note that I setup the "RestaurantReviews.Api" project to use
#    the local port 62770, and I expect this port to be used when you run
#    this solution. So in my below URLs I always use this 62770 local
#    port."
#
# is reworded to start with "Also note that ..." and gets a new trailing
# sentence about right-clicking the "RestaurantReviews.Api" project and
# choosing "Set as StartUp Project". (Everything else in the raw XML diff
# is Word's own spell/grammar-checker re-splitting runs and inserting
# <w:proofErr/> markers around already-present text - i.e. no visible
# content actually changes there.)

$d = $word.ActiveDocument

$leftDQ  = [char]0x201C   # “
$rightDQ = [char]0x201D   # ”

$old = "Finally note that I setup the " + $leftDQ + "RestaurantReviews.Api" + $rightDQ + `
    " project to use the local port 62770, and I expect this port to be used when you run this solution. So in my below URLs I always use this 62770 local port."

$new = "Also note that I setup the " + $leftDQ + "RestaurantReviews.Api" + $rightDQ + `
    " project to use the local port 62770, and I expect this port to be used when you run this solution. So in my below URLs I always use this 62770 local port." + `
    " Finally you will need to right-click on the " + $leftDQ + "RestaurantReviews.Api" + $rightDQ + `
    " project and select the " + $leftDQ + "Set as StartUp Project" + $rightDQ + " popup menu item. "

$replaced = $d.Content.Find.Execute($old, $true, $false, $false, $false, $false, $true, 1, $false, $new, 2)
Write-Output ("Reworded 'Finally note...' paragraph: " + $replaced)

# The document's hidden "_GoBack" bookmark (Word's "last edit position"
# marker) moves along with the edit: it used to sit in the later
# "reviewerEmail" paragraph and now sits right after the newly typed
# "popup menu item" text (just before the paragraph's closing ". ").
# Relocate it to match.
try {
    $oldMark = $d.Bookmarks("_GoBack")
    $oldMark.Delete()
    Write-Output "Removed old _GoBack bookmark"
} catch {
    Write-Output "No existing _GoBack bookmark found"
}

$tail = $d.Content
$foundTail = $tail.Find.Execute("popup menu item")
if ($foundTail) {
    $tail.Collapse(0)
    $d.Bookmarks.Add("_GoBack", $tail) | Out-Null
    Write-Output "Re-added _GoBack bookmark at new position"
}
